# "Permission for importing services passed"
#
# The services-import test fixture is updated so the sample row reflects a
# service (rather than advice) import whose referral method needs no
# special handling, and the "show referral disclaimer" flag is switched
# off now that the permission check passes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("type") for the sample row: advice -> service
$ws.Range("B2").Value = "service"

# Column R ("referral_method") for the sample row: internal -> none
$ws.Range("R2").Value = "none"

# Column Q ("show_referral_disclaimer") for the sample row: 1 -> 0
$ws.Range("Q2").Value = 0

# Leave the selection on the cell that was edited, like a user would.
$ws.Range("B2").Select()
